# Updated symbol list on Sat Jan 21 20:22:03 UTC 2023 with GitHub Actions
# Refresh the "Price" (D) and "Volume(1h)" (E) columns on the active
# (only) worksheet with the latest scraped values. These columns hold
# plain text (not numbers/percentages) in the source data, so we force
# a text number format before writing so Excel doesn't silently coerce
# "5.090" -> 5.09 or "2.36%" -> a percentage number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "304.09";  "E2"  = "2.36%"
    "D3"  = "35.59";   "E3"  = "12.68%"
    "D4"  = "5.090";   "E4"  = "2.42%"
    "D5"  = "0.07814"; "E5"  = "2.50%"
    "D6"  = "2.269";   "E6"  = "-0.10%"
    "D7"  = "8.121";   "E7"  = "3.60%"
    "D8"  = "4.020";   "E8"  = "6.18%"
    "D9"  = "0.9272";  "E9"  = "0.23%"
    "D10" = "0.09665"; "E10" = "0.59%"
    "D11" = "0.1823";  "E11" = "4.67%"
    "D12" = "0.08713"; "E12" = "3.62%"
    "D13" = "0.03419"; "E13" = "4.96%"
    "D14" = "0.09941"; "E14" = "0.03%"
    "D15" = "0.001480"; "E15" = "0.44%"
    "D16" = "0.005714"; "E16" = "0.60%"
    "D17" = "3.481";   "E17" = "-0.27%"
    "D18" = "2.147";   "E18" = "-2.23%"
    "D19" = "0.3418";  "E19" = "1.89%"
                        "E20" = "-0.08%"
    "D21" = "4.568";   "E21" = "12.61%"
    "D22" = "0.2236";  "E22" = "-2.00%"
    "D23" = "0.04679"; "E23" = "3.75%"
                        "E24" = "2.47%"
    "D25" = "0.004550"; "E25" = "5.07%"
    "D26" = "0.0001301"; "E26" = "1.08%"
    "D27" = "0.0002700"; "E27" = "-19.95%"
    "D39" = "0.01758"; "E39" = "5.94%"
    "D40" = "0.04699"; "E40" = "1.56%"
    "D41" = "0.007898"; "E41" = "5.67%"
    "D42" = "0.1421";  "E42" = "3.10%"
    "D43" = "0.008012"; "E43" = "-17.93%"
    "D44" = "0.002302"; "E44" = "7.66%"
    "D45" = "0.009098"; "E45" = "-3.49%"
    "D46" = "0.00006228"; "E46" = "2.50%"
                        "E47" = "0.60%"
    "D48" = "4.052";   "E48" = "52.64%"
                        "E49" = "35.32%"
                        "E50" = "0.60%"
                        "E51" = "0.60%"
}

# Make sure the target cells keep their existing plain-text storage
# (t="inlineStr"/"s") instead of being auto-converted to numbers or
# percentages by Excel's type inference.
$targetRange = $ws.Range("D2:E51")
$targetRange.NumberFormat = "@"

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# Restore the default (unstyled) look for the cells we touched so only
# the values themselves change, matching the original formatting.
$targetRange.Style = "Normal"
